$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 ---
# Seed formatting (border/fill/font = style index 3) by copying an existing data row.
$ws.Range("A7:F7").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

$ws.Range("A14").Value = 8

# B14: new date string -> plain text assignment is safe (not numeric-looking).
$ws.Range("B14").Value = "2021-02-11 10:54:03.0"

# C14: "12.12" already exists verbatim in C13 - copy it (value+type) to avoid
# Excel's text->number coercion that a fresh string assignment would trigger.
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4163)

# D14: new text "Egreso" -> safe as plain text.
$ws.Range("D14").Value = "Egreso"

# E14: same text already used in E13 - copy it across.
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4163)

# F14: "1" already exists verbatim in F7 - copy it across.
$ws.Range("F7").Copy()
$ws.Range("F14").PasteSpecial(-4163)

# --- Row 15 ---
$ws.Range("A7:F7").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

$ws.Range("A15").Value = 9

# B15: new date string.
$ws.Range("B15").Value = "2021-02-11 11:21:58.0"

# C15: "200" already exists verbatim in C8 - copy it across.
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# D15: "Ingreso" already exists verbatim in D7 - copy it across.
$ws.Range("D7").Copy()
$ws.Range("D15").PasteSpecial(-4163)

# E15: same text already used in E13.
$ws.Range("E13").Copy()
$ws.Range("E15").PasteSpecial(-4163)

# F15: "1" already exists verbatim in F7.
$ws.Range("F7").Copy()
$ws.Range("F15").PasteSpecial(-4163)
